# Update the cryptocurrency price/volume table with the latest scrape.
# Rows 2-33, 36-46: in-place price (D) / 1h-volume (E) refresh.
# Rows 34-35 and 47-51: coin ranking shuffled, so Coin/Link/Price/Volume are rewritten.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells are forced to Text format before assignment so Excel
# doesn't auto-convert numeric-looking strings (e.g. "303.40", "0.0795")
# into floating point numbers, which would lose the original formatting.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.891.58'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.229.26'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.40'
$ws.Range('E5').Value = '  -4.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.28'
$ws.Range('E6').Value = '  -7.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.565'
$ws.Range('E7').Value = '  -1.88%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.516'
$ws.Range('E9').Value = '  -7.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.80'
$ws.Range('E10').Value = '  -8.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0795'
$ws.Range('E11').Value = '  -4.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.05'
$ws.Range('E12').Value = '  -8.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.102'
$ws.Range('E13').Value = '  -3.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.568.84'
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.262.71'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.807'
$ws.Range('E16').Value = '  -5.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.40'
$ws.Range('E17').Value = '  -5.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.661.67'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0949'
$ws.Range('E19').Value = '  -3.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.10'
$ws.Range('E20').Value = '  -9.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.11'
$ws.Range('E21').Value = '  -6.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.05'
$ws.Range('E22').Value = '  -2.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.39'
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.88'
$ws.Range('E24').Value = '  -7.53%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.90'
$ws.Range('E26').Value = '  -10.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.71'
$ws.Range('E27').Value = '  -3.76%  '
$ws.Range('E28').Value = '  -1.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '35.87'
$ws.Range('E29').Value = '  -3.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.87'
$ws.Range('E30').Value = '  -5.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.70'
$ws.Range('E31').Value = '  -2.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.30'
$ws.Range('E32').Value = '  -4.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0794'
$ws.Range('E33').Value = '  -6.45%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.60'
$ws.Range('E34').Value = '  -3.45%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.21'
$ws.Range('E35').Value = '  +4.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.117'
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('E37').Value = '  -9.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.71'
$ws.Range('E38').Value = '  -11.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.32'
$ws.Range('E39').Value = '  -9.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.75'
$ws.Range('E40').Value = '  -11.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0296'
$ws.Range('E41').Value = '  -6.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.22'
$ws.Range('E42').Value = '  -13.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.01'
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.729.66'
$ws.Range('E44').Value = '  -3.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '82.98'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.184'
$ws.Range('E46').Value = '  -7.27%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '98.89'
$ws.Range('E47').Value = '  -4.47%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.90'
$ws.Range('E48').Value = '  -5.57%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.05'
$ws.Range('E49').Value = '  -3.74%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '67.24'
$ws.Range('E50').Value = '  -10.76%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '53.35'
$ws.Range('E51').Value = '  -8.58%  '
